# Apply the "EID mandatory resource" settings update to the Settings sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# 1) Update the description for "BirthdayList_Index_Birthday" row (C11) to the
#    new "EID" column-index description. (Shared string is registered first,
#    matching the order new strings were introduced in the saved file.)
$ws.Range("C11").Value = "Column index of ""EID"" in BirthdayList datatable"

# 2) Insert a new row above the old (empty) row 17, shifting the Email_Subject /
#    Email_Body rows (and everything below) down by one, and fill it with the
#    new "UpcomingBirthday_MandatoryEIDGetAll" setting.
$ws.Rows("17:17").Insert()
$ws.Range("B17").Value = "lorelie.a.pangan;sarah.c.c.layug;johnell.m.hernandez"
$ws.Range("A17").Value = "UpcomingBirthday_MandatoryEIDGetAll"
$ws.Range("C17").Value = "If one of the upcoming bday celebrant is here, all resources will be sent an email to ask for birthday greetings"

# Keep the new row's height consistent with the rest of the sheet.
$ws.Rows("17:17").RowHeight = 14.25

# 3) Reset the sheet selection back to the top-left cell (the previous "B20"
#    selection no longer makes sense after the insert).
$ws.Range("A1").Select()
